$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha 1")

# --- Update the column W header to reflect the revised Q1 classification ---
$ws.Range("W1").Value = "Q1 Type of ATD (Based on Li et al 2015)"

# --- Populate column W ("Q1 Type of ATD (Based on Li et al 2015)") per row ---
$w = @{
  2  = "Other"
  3  = "Complex architectural behavioral dependencies"
  4  = "Architecture Smells"
  5  = "Complex architectural behavioral dependencies"
  6  = "Complex architectural behavioral dependencies"
  7  = "System-level structure quality issues"
  8  = "Architectural compliance issues"
  9  = "Architectural compliance issues"
  10 = "Architecture Smells"
  11 = "Other"
  12 = "NA"
  14 = "NA"
  15 = "Architectural compliance issues"
  16 = "Architectural compliance issues"
  17 = "System-level structure quality issues"
  20 = "System-level structure quality issues"
  21 = "Architectural compliance issues"
  22 = "Architecture Smells"
  23 = "Architectural compliance issues"
  24 = "NA"
  25 = "System-level structure quality issues"
  26 = "Violations of good architectural practices"
  27 = "NA"
  28 = "System-level structure quality issues"
  29 = "Architecture Smells"
  30 = "Architecture Smells"
  32 = "Architectural compliance issues"
  33 = "NA"
  34 = "NA"
  35 = "NA"
  36 = "NA"
  39 = "NA"
  40 = "NA"
  41 = "Violations of good architectural practices"
  42 = "Violations of good architectural practices"
  44 = "System-level structure quality issues"
  45 = "Other"
  46 = "NA"
  47 = "NA"
}

foreach ($row in $w.Keys) {
    $ws.Range("W$row").Value = $w[$row]
}

# --- AB14 gets a "No" answer ---
$ws.Range("AB14").Value = "No"

# --- View/window/zoom/pane adjustments recorded by the author ---
$excel.ActiveWindow.Zoom = 150
$ws.Activate()
$ws.Range("U2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("W1").Select()

$wb.Windows.Item(1).WindowState = -4143
